# translations.xlsx update (improves #2059, #2061, #2063, #2054)
#
# 1. Adds a new "author" / "Autor" translation key, inserted as a new
#    row 13 (right after "date"/"Datum" and before "authors"/"Autoren"),
#    which pushes every row below it down by one.
# 2. Adds a new "reset" / "zurücksetzen" translation key as a new row
#    at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "author" row just above the existing "authors" row ---
$ws.Rows.Item(13).Insert()
$ws.Range("A13").Value2 = "author"
$ws.Range("B13").Value2 = "Autor"

# --- 2. Append the new "reset" row after the current last row ---
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
$newRow = $lastRow + 1
$ws.Range("A" + $newRow).Value2 = "reset"
$ws.Range("B" + $newRow).Value2 = "zurücksetzen"

# --- restore the view's active cell to sit near the newly inserted row ---
$ws.Range("B13").Select()
